$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to remain stored as text so values like "1.001" are not
# coerced into numbers by Excel's automatic type detection.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '30.350.75'
$ws.Range("E2").Value = '  -2.57%  '
$ws.Range("D3").Value = '1.942.01'
$ws.Range("E3").Value = '  -2.32%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.40%  '
$ws.Range("D5").Value = '251.28'
$ws.Range("E5").Value = '  -1.41%  '
$ws.Range("D6").Value = '0.7202'
$ws.Range("E6").Value = '  -10.07%  '
$ws.Range("E7").Value = '  +0.34%  '
$ws.Range("D8").Value = '0.3334'
$ws.Range("E8").Value = '  -4.74%  '
$ws.Range("D9").Value = '28.80'
$ws.Range("E9").Value = '  +2.34%  '
$ws.Range("D10").Value = '0.07352'
$ws.Range("E10").Value = '  +5.23%  '
$ws.Range("D11").Value = '0.8156'
$ws.Range("E11").Value = '  -3.56%  '
$ws.Range("D12").Value = '0.08133'
$ws.Range("E12").Value = '  -0.52%  '
$ws.Range("D13").Value = '1.939.87'
$ws.Range("E13").Value = '  -2.42%  '
$ws.Range("D14").Value = '5.494'
$ws.Range("E14").Value = '  -2.23%  '
$ws.Range("D15").Value = '95.17'
$ws.Range("E15").Value = '  -5.23%  '
$ws.Range("D16").Value = '14.94'
$ws.Range("E16").Value = '  -3.14%  '
$ws.Range("D17").Value = '0.000008417'
$ws.Range("E17").Value = '  +6.11%  '
$ws.Range("D18").Value = '30.345.53'
$ws.Range("E18").Value = '  -2.60%  '
$ws.Range("D19").Value = '253.10'
$ws.Range("E19").Value = '  -7.31%  '
$ws.Range("D20").Value = '5.897'
$ws.Range("E20").Value = '  +0.46%  '
$ws.Range("D21").Value = '2.195.01'
$ws.Range("E21").Value = '  -2.43%  '
$ws.Range("E22").Value = '  +0.29%  '
$ws.Range("D23").Value = '1.003'
$ws.Range("E23").Value = '  +0.55%  '
$ws.Range("D24").Value = '6.981'
$ws.Range("E24").Value = '  -1.00%  '
$ws.Range("D25").Value = '9.850'
$ws.Range("E25").Value = '  -1.39%  '
$ws.Range("D26").Value = '163.05'
$ws.Range("E26").Value = '  -1.53%  '
$ws.Range("D27").Value = '2.413'
$ws.Range("E27").Value = '  +2.96%  '
$ws.Range("D28").Value = '19.39'
$ws.Range("E28").Value = '  -2.57%  '
$ws.Range("D29").Value = '0.1322'
$ws.Range("E29").Value = '  -12.66%  '
$ws.Range("E30").Value = '  -1.47%  '
$ws.Range("D31").Value = '1.346'
$ws.Range("E31").Value = '  -0.67%  '
$ws.Range("D32").Value = '4.460'
$ws.Range("E32").Value = '  -2.67%  '
$ws.Range("D33").Value = '4.255'
$ws.Range("E33").Value = '  -3.50%  '
$ws.Range("D34").Value = '0.05270'
$ws.Range("E34").Value = '  +0.45%  '
$ws.Range("D35").Value = '1.302'
$ws.Range("E35").Value = '  +6.77%  '
$ws.Range("D36").Value = '0.7572'
$ws.Range("E36").Value = '  -2.82%  '
$ws.Range("D37").Value = '2.747'
$ws.Range("E37").Value = '  -0.56%  '
$ws.Range("D38").Value = '0.01995'
$ws.Range("E38").Value = '  -0.18%  '
$ws.Range("D39").Value = '2.857'
$ws.Range("E39").Value = '  -0.83%  '
$ws.Range("D40").Value = '81.05'
$ws.Range("E40").Value = '  +1.65%  '
$ws.Range("D41").Value = '6.622'
$ws.Range("E41").Value = '  -0.47%  '
$ws.Range("E42").Value = '  -2.36%  '
$ws.Range("D43").Value = '2.041'
$ws.Range("E43").Value = '  -4.07%  '
$ws.Range("D44").Value = '0.8476'
$ws.Range("E44").Value = '  -0.59%  '
$ws.Range("E45").Value = '  +0.36%  '
$ws.Range("D46").Value = '102.70'
$ws.Range("E46").Value = '  -1.92%  '
$ws.Range("D47").Value = '9.817'
$ws.Range("E47").Value = '  -0.47%  '
$ws.Range("D48").Value = '7.502'
$ws.Range("E48").Value = '  -2.23%  '
$ws.Range("D49").Value = '36.91'
$ws.Range("E49").Value = '  +0.34%  '
$ws.Range("D50").Value = '0.4194'
$ws.Range("E50").Value = '  -2.35%  '
$ws.Range("D51").Value = '1.506'
$ws.Range("E51").Value = '  -1.88%  '
